$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "떡"
$ws.Range("B3").Value = "간장"
$ws.Range("B4").Value = "식초"
$ws.Range("B5").Value = "간장"
$ws.Range("B6").Value = "삽겹살"
$ws.Range("B7").Value = "꿀"
